$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = '''27.678.13'
$ws.Range("E2").Value2 = '  +1.43%  '
$ws.Range("D3").Value2 = '''1.864.84'
$ws.Range("E3").Value2 = '  +0.50%  '
$ws.Range("D4").Value2 = '''1.005'
$ws.Range("E4").Value2 = '  +0.51%  '
$ws.Range("D5").Value2 = '''331.66'
$ws.Range("E5").Value2 = '  +3.34%  '
$ws.Range("D6").Value2 = '''1.005'
$ws.Range("E6").Value2 = '  +0.45%  '
$ws.Range("D8").Value2 = '''0.3931'
$ws.Range("E8").Value2 = '  +2.10%  '
$ws.Range("D9").Value2 = '''49.01'
$ws.Range("E9").Value2 = '  +1.72%  '
$ws.Range("D10").Value2 = '''0.07987'
$ws.Range("E10").Value2 = '  +1.30%  '
$ws.Range("D11").Value2 = '''1.020'
$ws.Range("E11").Value2 = '  +0.50%  '
$ws.Range("E12").Value2 = '  +1.88%  '
$ws.Range("D13").Value2 = '''1.866.29'
$ws.Range("E13").Value2 = '  +0.63%  '
$ws.Range("D14").Value2 = '''5.924'
$ws.Range("E14").Value2 = '  +1.07%  '
$ws.Range("D15").Value2 = '''7.113'
$ws.Range("E15").Value2 = '  -0.47%  '
$ws.Range("D16").Value2 = '''1.008'
$ws.Range("E16").Value2 = '  +0.76%  '
$ws.Range("E17").Value2 = '  +1.83%  '
$ws.Range("D18").Value2 = '''86.60'
$ws.Range("E18").Value2 = '  +1.08%  '
$ws.Range("D19").Value2 = '''0.06654'
$ws.Range("E19").Value2 = '  +2.22%  '
$ws.Range("D20").Value2 = '''17.20'
$ws.Range("E20").Value2 = '  +1.31%  '
$ws.Range("D21").Value2 = '''1.005'
$ws.Range("E21").Value2 = '  +0.47%  '
$ws.Range("D22").Value2 = '''27.680.39'
$ws.Range("E22").Value2 = '  +1.40%  '
$ws.Range("D23").Value2 = '''5.481'
$ws.Range("E23").Value2 = '  -0.04%  '
$ws.Range("D24").Value2 = '''10.97'
$ws.Range("E24").Value2 = '  +1.77%  '
$ws.Range("D25").Value2 = '''2.310'
$ws.Range("E25").Value2 = '  +2.17%  '
$ws.Range("D26").Value2 = '''2.071.26'
$ws.Range("E26").Value2 = '  -0.23%  '
$ws.Range("D27").Value2 = '''159.02'
$ws.Range("E27").Value2 = '  +5.06%  '
$ws.Range("D28").Value2 = '''20.14'
$ws.Range("E28").Value2 = '  +0.91%  '
$ws.Range("D29").Value2 = '''2.085'
$ws.Range("E29").Value2 = '  +1.45%  '
$ws.Range("D30").Value2 = '''5.538'
$ws.Range("E30").Value2 = '  +0.36%  '
$ws.Range("D31").Value2 = '''123.44'
$ws.Range("E31").Value2 = '  +2.59%  '
$ws.Range("D32").Value2 = '''0.9654'
$ws.Range("E32").Value2 = '  +2.73%  '
$ws.Range("D33").Value2 = '''0.09499'
$ws.Range("E33").Value2 = '  +2.35%  '
$ws.Range("D34").Value2 = '''1.442'
$ws.Range("E34").Value2 = '  -1.64%  '
$ws.Range("D35").Value2 = '''3.602'
$ws.Range("E35").Value2 = '  +1.04%  '
$ws.Range("E36").Value2 = '  -0.14%  '
$ws.Range("D37").Value2 = '''0.02249'
$ws.Range("E37").Value2 = '  +1.32%  '
$ws.Range("D38").Value2 = '''0.06055'
$ws.Range("E38").Value2 = '  +1.21%  '
$ws.Range("D39").Value2 = '''1.228'
$ws.Range("E39").Value2 = '  +2.57%  '
$ws.Range("D40").Value2 = '''8.102'
$ws.Range("E40").Value2 = '  -2.53%  '
$ws.Range("D41").Value2 = '''1.004'
$ws.Range("E41").Value2 = '  +0.47%  '
$ws.Range("D42").Value2 = '''0.5968'
$ws.Range("E42").Value2 = '  +1.05%  '
$ws.Range("D43").Value2 = '''0.1890'
$ws.Range("E43").Value2 = '  +0.84%  '
$ws.Range("D44").Value2 = '''10.21'
$ws.Range("E44").Value2 = '  +0.65%  '
$ws.Range("D45").Value2 = '''1.271'
$ws.Range("E45").Value2 = '  +0.33%  '
$ws.Range("D46").Value2 = '''0.5685'
$ws.Range("E46").Value2 = '  +0.98%  '
$ws.Range("D47").Value2 = '''12.25'
$ws.Range("E47").Value2 = '  +2.60%  '
$ws.Range("D48").Value2 = '''3.395'
$ws.Range("E48").Value2 = '  +1.23%  '
$ws.Range("D49").Value2 = '''1.933'
$ws.Range("E49").Value2 = '  +0.63%  '
$ws.Range("D50").Value2 = '''0.06827'
$ws.Range("E50").Value2 = '  +0.10%  '
$ws.Range("D51").Value2 = '''113.50'
$ws.Range("E51").Value2 = '  +4.60%  '
